$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.987.39"
$ws.Range("D3").Value = "'3.448.85"
$ws.Range("E3").Value = "'  -1.46%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'577.44"
$ws.Range("E5").Value = "'  -1.30%  "
$ws.Range("D6").Value = "'148.44"
$ws.Range("E6").Value = "'  +0.89%  "
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("D9").Value = "'7.81"
$ws.Range("E9").Value = "'  +1.67%  "
$ws.Range("E10").Value = "'  -2.65%  "
$ws.Range("E11").Value = "'  +1.70%  "
$ws.Range("D12").Value = "'4.040.22"
$ws.Range("E12").Value = "'  -1.47%  "
$ws.Range("E13").Value = "'  +2.60%  "
$ws.Range("D14").Value = "'28.54"
$ws.Range("E14").Value = "'  -3.70%  "
$ws.Range("D15").Value = "'3.444.81"
$ws.Range("E15").Value = "'  -1.74%  "
$ws.Range("E16").Value = "'  -2.29%  "
$ws.Range("D17").Value = "'63.003.48"
$ws.Range("E17").Value = "'  -0.59%  "
$ws.Range("D18").Value = "'6.39"
$ws.Range("E18").Value = "'  +1.68%  "
$ws.Range("D19").Value = "'14.33"
$ws.Range("E19").Value = "'  -0.04%  "
$ws.Range("D20").Value = "'9.11"
$ws.Range("E20").Value = "'  -3.69%  "
$ws.Range("D21").Value = "'384.91"
$ws.Range("E21").Value = "'  -2.59%  "
$ws.Range("D22").Value = "'0.558"
$ws.Range("E22").Value = "'  -1.58%  "
$ws.Range("D23").Value = "'74.39"
$ws.Range("E23").Value = "'  -1.35%  "
$ws.Range("E24").Value = "'  -0.13%  "
$ws.Range("D25").Value = "'3.582.96"
$ws.Range("E25").Value = "'  -1.63%  "
$ws.Range("E26").Value = "'  -5.43%  "
$ws.Range("E27").Value = "'  -5.30%  "
$ws.Range("D28").Value = "'7.66"
$ws.Range("E28").Value = "'  -2.12%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "'  +1.77%  "
$ws.Range("D30").Value = "'7.99"
$ws.Range("E30").Value = "'  -3.47%  "
$ws.Range("E31").Value = "'  -3.01%  "
$ws.Range("E32").Value = "'  +0.02%  "
$ws.Range("E33").Value = "'  -2.75%  "
$ws.Range("D34").Value = "'1.30"
$ws.Range("E34").Value = "'  -9.18%  "
$ws.Range("D35").Value = "'5.36"
$ws.Range("E35").Value = "'  +0.09%  "
$ws.Range("D36").Value = "'1.61"
$ws.Range("E36").Value = "'  +1.56%  "
$ws.Range("B37").Value = "'EnergySwap"
$ws.Range("C37").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "'31.92"
$ws.Range("E37").Value = "'  -2.75%  "
$ws.Range("B38").Value = "'Aptos"
$ws.Range("C38").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'7.03"
$ws.Range("E38").Value = "'  -2.62%  "
$ws.Range("D39").Value = "'169.53"
$ws.Range("E39").Value = "'  -2.19%  "
$ws.Range("D40").Value = "'3.483.14"
$ws.Range("E40").Value = "'  -1.48%  "
$ws.Range("D41").Value = "'0.0763"
$ws.Range("E41").Value = "'  -1.46%  "
$ws.Range("E42").Value = "'  -1.59%  "
$ws.Range("D43").Value = "'42.46"
$ws.Range("E43").Value = "'  -0.20%  "
$ws.Range("E44").Value = "'  -1.81%  "
$ws.Range("D45").Value = "'4.35"
$ws.Range("E45").Value = "'  -4.04%  "
$ws.Range("E46").Value = "'  -1.92%  "
$ws.Range("D47").Value = "'2.574.67"
$ws.Range("E47").Value = "'  -0.77%  "
$ws.Range("D48").Value = "'2.28"
$ws.Range("E48").Value = "'  +0.68%  "
$ws.Range("E49").Value = "'  +1.11%  "
$ws.Range("B50").Value = "'InjectiveProtocol"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'22.48"
$ws.Range("E50").Value = "'  -6.34%  "
$ws.Range("B51").Value = "'FirstDigitalUSD"
$ws.Range("C51").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "'  -0.04%  "
